# Slide 16 ("PP") + ": Risks" title: merge the two runs into a single
# run reading "PP: Risks", dropping the now-redundant endParaRPr.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(16)
$tr = $s.Shapes.Item(1).TextFrame.TextRange
$tr.Delete()
$tr.Text = "PP: Risks"
